# Add chapter 8.4 -> lab 13 (column K) and chapter 7.5 -> lab 12 (column J),
# and remove the row for chapter "8.6" entirely (it no longer exists after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 holds chapter "8.6" (A14) - delete the whole row, shifting everything
# below it up by one (this also fixes up the shared formulas, shared-string
# indices, dimension, and the SUM/COUNT average formula's dependents).
$ws.Rows("14:14").Delete() | Out-Null

# Chapter "7.5" is row 10 (A10) - mark it done for lab 12 (column J).
$ws.Range("J10").Value = 1

# Chapter "8.4" is row 13 (A13) - mark it done for lab 13 (column K).
$ws.Range("K13").Value = 1

# Shrink the conditional formatting range that tracked N3:N37 down to N3:N36
# to match the now-shorter table.
$cf = $ws.Range("N3:N1000").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("N3:N36"))

# Match the author's final selection (cell J14, which is now the row for
# chapter "9.2" after the row-14 deletion).
$ws.Range("J14").Select() | Out-Null
